$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pension recipient counts for 2015-2018 (columns E-H, row 4)
$ws.Range("E4").Value = 27784
$ws.Range("F4").Value = 28447
$ws.Range("G4").Value = 29070
$ws.Range("H4").Value = 29523

# Move / extend the active selection to E4:H4 (active cell E4)
$ws.Range("E4:H4").Select()
